# Auto-generated edit script: refresh computed market-price / profit columns (H..N)
# across multiple Leve sheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 771681.3
$ws.Range("I19").Value = 2001295.8
$ws.Range("J19").Value = 3172.25
$ws.Range("K19").Value = 2001295.8
$ws.Range("L19").Value = 3172.25
$ws.Range("M19").Value = -2001120.8
$ws.Range("N19").Value = -3522.25
$ws.Range("H43").Value = 1587.9286
$ws.Range("J43").Value = 1499.75
$ws.Range("L43").Value = 1499.75
$ws.Range("N43").Value = -1637.75
$ws.Range("H98").Value = 1047.1316
$ws.Range("I98").Value = 953.8108
$ws.Range("K98").Value = 953.8108
$ws.Range("M98").Value = 544.1892
$ws.Range("H107").Value = 375.18182
$ws.Range("I107").Value = 402.3
$ws.Range("K107").Value = 402.3
$ws.Range("M107").Value = 1517.7
$ws.Range("H112").Value = 1920.9412
$ws.Range("J112").Value = 1947.36
$ws.Range("L112").Value = 5842.08
$ws.Range("N112").Value = -8058.08
$ws.Range("H122").Value = 1047.1316
$ws.Range("I122").Value = 953.8108
$ws.Range("K122").Value = 2861.4324
$ws.Range("M122").Value = -411.4323999999997
$ws.Range("H131").Value = 1628.75
$ws.Range("J131").Value = 2869.5
$ws.Range("L131").Value = 8608.5
$ws.Range("N131").Value = -18688.5
$ws.Range("H137").Value = 73443.21000000001
$ws.Range("I137").Value = 1719.4
$ws.Range("J137").Value = 113289.78
$ws.Range("K137").Value = 5158.200000000001
$ws.Range("L137").Value = 339869.34
$ws.Range("M137").Value = -2608.200000000001
$ws.Range("N137").Value = -344969.34
$ws.Range("H138").Value = 3489.4
$ws.Range("I138").Value = 4698.8
$ws.Range("J138").Value = 2971.0857
$ws.Range("K138").Value = 14096.4
$ws.Range("L138").Value = 8913.257100000001
$ws.Range("M138").Value = -8956.400000000001
$ws.Range("N138").Value = -19193.2571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7053.54
$ws.Range("I32").Value = 6214.274
$ws.Range("J32").Value = 22999.6
$ws.Range("K32").Value = 6214.274
$ws.Range("L32").Value = 22999.6
$ws.Range("M32").Value = -5927.274
$ws.Range("N32").Value = -23573.6
$ws.Range("H59").Value = 12300
$ws.Range("I59").Value = 8600
$ws.Range("J59").Value = 16000
$ws.Range("K59").Value = 8600
$ws.Range("L59").Value = 16000
$ws.Range("M59").Value = -7796
$ws.Range("N59").Value = -17608
$ws.Range("H122").Value = 18613.773
$ws.Range("I122").Value = 21242.316
$ws.Range("K122").Value = 63726.948
$ws.Range("M122").Value = -61276.948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 251529.88
$ws.Range("I86").Value = 1748.4286
$ws.Range("K86").Value = 1748.4286
$ws.Range("M86").Value = -625.4286
$ws.Range("H89").Value = 251529.88
$ws.Range("I89").Value = 1748.4286
$ws.Range("K89").Value = 8742.143
$ws.Range("M89").Value = -3126.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2673.64
$ws.Range("I31").Value = 2366.5557
$ws.Range("J31").Value = 3463.2856
$ws.Range("K31").Value = 2366.5557
$ws.Range("L31").Value = 3463.2856
$ws.Range("M31").Value = -2071.5557
$ws.Range("N31").Value = -4053.2856
$ws.Range("H34").Value = 2673.64
$ws.Range("I34").Value = 2366.5557
$ws.Range("J34").Value = 3463.2856
$ws.Range("K34").Value = 2366.5557
$ws.Range("L34").Value = 3463.2856
$ws.Range("M34").Value = -2164.5557
$ws.Range("N34").Value = -3867.2856
$ws.Range("H41").Value = 29000
$ws.Range("J41").Value = 29000
$ws.Range("L41").Value = 29000
$ws.Range("N41").Value = -29856
$ws.Range("H59").Value = 17971.428
$ws.Range("J59").Value = 17971.428
$ws.Range("L59").Value = 17971.428
$ws.Range("N59").Value = -20261.428
$ws.Range("H132").Value = 1455.7037
$ws.Range("I132").Value = 1190.0952
$ws.Range("J132").Value = 2385.3333
$ws.Range("K132").Value = 3570.2856
$ws.Range("L132").Value = 7155.999899999999
$ws.Range("M132").Value = -1040.2856
$ws.Range("N132").Value = -12215.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 117.8
$ws.Range("J2").Value = 95
$ws.Range("L2").Value = 570
$ws.Range("N2").Value = -796
$ws.Range("H17").Value = 13401.2
$ws.Range("J17").Value = 22002
$ws.Range("L17").Value = 66006
$ws.Range("N17").Value = -66344
$ws.Range("H23").Value = 257.66666
$ws.Range("J23").Value = 292.5
$ws.Range("L23").Value = 877.5
$ws.Range("N23").Value = -1347.5
$ws.Range("H33").Value = 254.09525
$ws.Range("J33").Value = 323.18182
$ws.Range("L33").Value = 1939.09092
$ws.Range("N33").Value = -2505.09092
$ws.Range("H37").Value = 54666.668
$ws.Range("J37").Value = 54666.668
$ws.Range("L37").Value = 164000.004
$ws.Range("N37").Value = -164224.004
$ws.Range("H56").Value = 6639.4165
$ws.Range("I56").Value = 6639.4165
$ws.Range("K56").Value = 6639.4165
$ws.Range("M56").Value = -6109.4165
$ws.Range("H116").Value = 2426.7144
$ws.Range("I116").Value = 1335.8
$ws.Range("J116").Value = 3032.7778
$ws.Range("K116").Value = 4007.4
$ws.Range("L116").Value = 9098.3334
$ws.Range("M116").Value = -565.3999999999996
$ws.Range("N116").Value = -15982.3334
$ws.Range("H122").Value = 1486.2667
$ws.Range("J122").Value = 1724.5
$ws.Range("L122").Value = 15520.5
$ws.Range("N122").Value = -20420.5
$ws.Range("H131").Value = 12844524
$ws.Range("I131").Value = 71428984
$ws.Range("K131").Value = 214286952
$ws.Range("M131").Value = -214281912
$ws.Range("H132").Value = 6189.7
$ws.Range("J132").Value = 7507.5
$ws.Range("L132").Value = 67567.5
$ws.Range("N132").Value = -72627.5
$ws.Range("H136").Value = 1273.5883
$ws.Range("I136").Value = 1273.5883
$ws.Range("K136").Value = 3820.7649
$ws.Range("M136").Value = 1279.2351
$ws.Range("H137").Value = 6775.4707
$ws.Range("J137").Value = 7174.5
$ws.Range("L137").Value = 21523.5
$ws.Range("N137").Value = -31723.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 35000
$ws.Range("J62").Value = 35000
$ws.Range("L62").Value = 35000
$ws.Range("N62").Value = -36372
$ws.Range("H65").Value = 35000
$ws.Range("J65").Value = 35000
$ws.Range("L65").Value = 105000
$ws.Range("N65").Value = -111864
$ws.Range("H97").Value = 2786
$ws.Range("I97").Value = 2905.111
$ws.Range("J97").Value = 2250
$ws.Range("K97").Value = 2905.111
$ws.Range("L97").Value = 2250
$ws.Range("M97").Value = -2409.111
$ws.Range("N97").Value = -3242
$ws.Range("H102").Value = 1380.9354
$ws.Range("I102").Value = 1360.3
$ws.Range("K102").Value = 1360.3
$ws.Range("M102").Value = 261.7
$ws.Range("H113").Value = 1749.75
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1999.6666
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1999.6666
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -6339.6666
$ws.Range("H122").Value = 1176.5454
$ws.Range("I122").Value = 429.2857
$ws.Range("J122").Value = 1525.2667
$ws.Range("K122").Value = 1287.8571
$ws.Range("L122").Value = 4575.800099999999
$ws.Range("M122").Value = 1162.1429
$ws.Range("N122").Value = -9475.8001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 17185
$ws.Range("I40").Value = 20437.875
$ws.Range("K40").Value = 20437.875
$ws.Range("M40").Value = -20301.875
$ws.Range("H61").Value = 2701.6155
$ws.Range("I61").Value = 2514.5557
$ws.Range("K61").Value = 2514.5557
$ws.Range("M61").Value = -2312.5557
$ws.Range("H113").Value = 2701.6155
$ws.Range("I113").Value = 2514.5557
$ws.Range("K113").Value = 2514.5557
$ws.Range("M113").Value = -344.5556999999999
$ws.Range("H122").Value = 3200.2307
$ws.Range("J122").Value = 3649.8333
$ws.Range("L122").Value = 10949.4999
$ws.Range("N122").Value = -15849.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H122").Value = 30605.072
$ws.Range("I122").Value = 46080.168
$ws.Range("J122").Value = 2749.9
$ws.Range("K122").Value = 138240.504
$ws.Range("L122").Value = 8249.700000000001
$ws.Range("M122").Value = -135790.504
$ws.Range("N122").Value = -13149.7
